$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task list (rows 2-18, column A) ---
$tasks = @(
    "Normalization",
    "Fill Missing Values",
    "Outlier Removal",
    "Numeric to Category",
    "Principal Component Analysis",
    "Change Column Type",
    "Rename Column",
    "Special Character Removal",
    "Trim Whitespace",
    "Replace Substrings",
    "Text Case",
    "Remove Stopwords",
    "Collapse Rare Categories",
    "Tokenization",
    "Regex",
    "Datetime Components",
    "Remove Columns"
)

# --- Header row ---
$ws.Range("A1").Value = "Task"
$ws.Range("B1").Value = "Local IDE ( PYTHON ) in minutes"
$ws.Range("C1").Value = "Data Polish"

# --- Body rows (2..18) ---
for ($i = 0; $i -lt $tasks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tasks[$i]
}

# --- Footer rows ---
$ws.Range("A19").Value = "TOTAL"
$ws.Range("A20").Value = "TOTAL in hours"

# --- Fonts / styles ---
# Every used cell is 16pt Calibri; set this across the whole range first so
# only a single (non-bold) 16pt font gets created.
$ws.Range("A1:C20").Font.Size = 16

# Header row + TOTAL rows (col A & B only) are additionally bold -- this
# derives a second (bold) font from the already-16pt cells.
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A19:B19").Font.Bold = $true
$ws.Range("A20:B20").Font.Bold = $true

# --- Row heights ---
$ws.Range("A1:C20").RowHeight = 21

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 33.83203125
$ws.Columns.Item(2).ColumnWidth = 38

# --- Selection ---
$ws.Range("B2:B20").Select
